# Auto-generated edit script
# Re-applies the updated FFXIV crafting-leve market/profit figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets,
# matching the refreshed scheduled-runner market data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 347.5
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H113").Value = 2644.25
$ws.Range("I113").Value = 2359
$ws.Range("K113").Value = 2359
$ws.Range("M113").Value = 895
$ws.Range("H116").Value = 998.3333
$ws.Range("I116").Value = 998.3333
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 998.3333
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 2443.6667
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2750
$ws.Range("I26").Value = 2750
$ws.Range("K26").Value = 2750
$ws.Range("M26").Value = -2420
$ws.Range("H35").Value = 900
$ws.Range("I35").Value = 900
$ws.Range("K35").Value = 900
$ws.Range("M35").Value = -494
$ws.Range("H61").Value = 4176.6
$ws.Range("I61").Value = 2514.9
$ws.Range("K61").Value = 2514.9
$ws.Range("M61").Value = -2302.9
$ws.Range("H132").Value = 7038.8
$ws.Range("I132").Value = 6798.5
$ws.Range("K132").Value = 20395.5
$ws.Range("M132").Value = -17865.5
$ws.Range("H136").Value = 4176.6
$ws.Range("I136").Value = 2514.9
$ws.Range("K136").Value = 7544.700000000001
$ws.Range("M136").Value = -4994.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 1587.5
$ws.Range("I33").Value = 1587.5
$ws.Range("K33").Value = 1587.5
$ws.Range("M33").Value = -1251.5
$ws.Range("H86").Value = 6134.615
$ws.Range("I86").Value = 3731.25
$ws.Range("K86").Value = 3731.25
$ws.Range("M86").Value = -2608.25
$ws.Range("H87").Value = 175000
$ws.Range("J87").Value = 175000
$ws.Range("L87").Value = 175000
$ws.Range("N87").Value = -177496
$ws.Range("H89").Value = 6134.615
$ws.Range("I89").Value = 3731.25
$ws.Range("K89").Value = 18656.25
$ws.Range("M89").Value = -13040.25
$ws.Range("H90").Value = 175000
$ws.Range("J90").Value = 175000
$ws.Range("L90").Value = 525000
$ws.Range("N90").Value = -537480
$ws.Range("H107").Value = 4160.2
$ws.Range("I107").Value = 2854.818
$ws.Range("K107").Value = 2854.818
$ws.Range("M107").Value = -934.8180000000002
$ws.Range("H134").Value = 3628.6667
$ws.Range("I134").Value = 3628.6667
$ws.Range("K134").Value = 10886.0001
$ws.Range("M134").Value = -8351.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 1215
$ws.Range("I33").Value = 1215
$ws.Range("K33").Value = 1215
$ws.Range("M33").Value = -836
$ws.Range("H122").Value = 1859.1818
$ws.Range("I122").Value = 1751.5714
$ws.Range("J122").Value = 2047.5
$ws.Range("K122").Value = 5254.7142
$ws.Range("L122").Value = 6142.5
$ws.Range("M122").Value = -2804.7142
$ws.Range("N122").Value = -11042.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 9999
$ws.Range("J76").Value = 9999
$ws.Range("L76").Value = 29997
$ws.Range("N76").Value = -30763
$ws.Range("H79").Value = 9999
$ws.Range("J79").Value = 9999
$ws.Range("L79").Value = 29997
$ws.Range("N79").Value = -32649
$ws.Range("H81").Value = 1869.8
$ws.Range("I81").Value = 850
$ws.Range("J81").Value = 2124.75
$ws.Range("K81").Value = 2550
$ws.Range("L81").Value = 6374.25
$ws.Range("M81").Value = -1427
$ws.Range("N81").Value = -8620.25
$ws.Range("H84").Value = 1869.8
$ws.Range("I84").Value = 850
$ws.Range("J84").Value = 2124.75
$ws.Range("K84").Value = 7650
$ws.Range("L84").Value = 19122.75
$ws.Range("M84").Value = -2034
$ws.Range("N84").Value = -30354.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 24800
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H69").Value = 65000
$ws.Range("J69").Value = 65000
$ws.Range("L69").Value = 65000
$ws.Range("N69").Value = -66498
$ws.Range("H72").Value = 65000
$ws.Range("J72").Value = 65000
$ws.Range("L72").Value = 195000
$ws.Range("N72").Value = -202488
$ws.Range("H94").Value = 27464.25
$ws.Range("J94").Value = 27464.25
$ws.Range("L94").Value = 27464.25
$ws.Range("N94").Value = -28816.25
$ws.Range("H102").Value = 3881.25
$ws.Range("I102").Value = 5525
$ws.Range("J102").Value = 3333.3333
$ws.Range("K102").Value = 5525
$ws.Range("L102").Value = 3333.3333
$ws.Range("M102").Value = -3903
$ws.Range("N102").Value = -6577.3333
$ws.Range("H108").Value = 50000
$ws.Range("I108").Value = 50000
$ws.Range("K108").Value = 50000
$ws.Range("M108").Value = -46160
$ws.Range("H132").Value = 5196.2
$ws.Range("I132").Value = 4394.5
$ws.Range("K132").Value = 13183.5
$ws.Range("M132").Value = -10653.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51622
$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -158112
$ws.Range("H132").Value = 3659.8
$ws.Range("I132").Value = 3659.8
$ws.Range("K132").Value = 10979.4
$ws.Range("M132").Value = -8449.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 9250
$ws.Range("I32").Value = 13300
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 13300
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = -12983
$ws.Range("N32").Value = -3134
$ws.Range("H52").Value = 23995
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H54").Value = 18499.834
$ws.Range("I54").Value = 2200
$ws.Range("K54").Value = 2200
$ws.Range("M54").Value = -1680
$ws.Range("H64").Value = 90000
$ws.Range("I64").Value = 90000
$ws.Range("K64").Value = 90000
$ws.Range("M64").Value = -89752
$ws.Range("H67").Value = 90000
$ws.Range("I67").Value = 90000
$ws.Range("K67").Value = 90000
$ws.Range("M67").Value = -89142
$ws.Range("H122").Value = 2842.4285
$ws.Range("I122").Value = 2166.6667
$ws.Range("K122").Value = 6500.000100000001
$ws.Range("M122").Value = -4050.000100000001
$ws.Range("H136").Value = 3714
$ws.Range("I136").Value = 3071
$ws.Range("K136").Value = 9213
$ws.Range("M136").Value = -6663
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

Write-Host "Applied Sheets updates via scheduled runner (190 cells)"
